$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated cryptos list (price + volume(1h) refresh, two coin rows swapped positions) ---

# Cells whose new text is a plain (non-numeric-looking) string: assign directly.
$ws.Range("D2").Value = '55.918.60'
$ws.Range("E2").Value = '  +8.58%  '
$ws.Range("D3").Value = '3.224.64'
$ws.Range("E3").Value = '  +4.03%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("E5").Value = '  +3.83%  '
$ws.Range("E6").Value = '  +6.19%  '
$ws.Range("E7").Value = '  +2.79%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +6.69%  '
$ws.Range("E10").Value = '  +6.25%  '
$ws.Range("E11").Value = '  +5.63%  '
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '3.726.12'
$ws.Range("E13").Value = '  +3.80%  '
$ws.Range("E14").Value = '  +2.85%  '
$ws.Range("E15").Value = '  +2.95%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.218.28'
$ws.Range("E16").Value = '  +3.85%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("E17").Value = '  +6.19%  '
$ws.Range("E18").Value = '  -4.22%  '
$ws.Range("D19").Value = '55.706.54'
$ws.Range("E19").Value = '  +8.05%  '
$ws.Range("E20").Value = '  +2.74%  '
$ws.Range("E21").Value = '  +6.22%  '
$ws.Range("E22").Value = '  +5.51%  '
$ws.Range("E23").Value = '  +13.90%  '
$ws.Range("E24").Value = '  +6.95%  '
$ws.Range("E25").Value = '  +2.43%  '
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("E27").Value = '  +5.04%  '
$ws.Range("E28").Value = '  +3.07%  '
$ws.Range("E29").Value = '  +2.29%  '
$ws.Range("E31").Value = '  +10.41%  '
$ws.Range("E33").Value = '  +4.36%  '
$ws.Range("E34").Value = '  +3.18%  '
$ws.Range("E35").Value = '  +3.29%  '
$ws.Range("E36").Value = '  +2.09%  '
$ws.Range("E37").Value = '  +5.39%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("E38").Value = '  +23.36%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E40").Value = '  +3.70%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E41").Value = '  +2.24%  '
$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E42").Value = '  +10.25%  '
$ws.Range("E43").Value = '  +3.50%  '
$ws.Range("E44").Value = '  +2.79%  '
$ws.Range("E45").Value = '  -4.71%  '
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E47").Value = '  +2.73%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.162.85'
$ws.Range("E48").Value = '  +4.83%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("E49").Value = '  +42.97%  '
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("E51").Value = '  +9.37%  '

# Cells whose new text looks like a plain number (e.g. "0.999", "39.29"): the source
# data keeps these as literal text, so force text format first to avoid Excel
# re-interpreting them as numeric values.
$numericTextCells = @("D4","D5","D6","D10","D11","D12","D15","D17","D18","D20","D21","D23","D24","D25","D26","D27","D28","D30","D31","D34","D37","D38","D39","D40","D41","D42","D46","D49","D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").Value = '400.15'
$ws.Range("D6").Value = '109.69'
$ws.Range("D10").Value = '39.29'
$ws.Range("D11").Value = '0.0903'
$ws.Range("D12").Value = '0.141'
$ws.Range("D15").Value = '8.07'
$ws.Range("D17").Value = '1.06'
$ws.Range("D18").Value = '10.63'
$ws.Range("D20").Value = '3.39'
$ws.Range("D21").Value = '0.0000102'
$ws.Range("D23").Value = '303.06'
$ws.Range("D24").Value = '74.78'
$ws.Range("D25").Value = '3.23'
$ws.Range("D26").Value = '8.23'
$ws.Range("D27").Value = '28.40'
$ws.Range("D28").Value = '7.51'
$ws.Range("D30").Value = '0.999'
$ws.Range("D31").Value = '11.42'
$ws.Range("D34").Value = '36.41'
$ws.Range("D37").Value = '3.54'
$ws.Range("D38").Value = '3.10'
$ws.Range("D39").Value = '0.998'
$ws.Range("D40").Value = '133.72'
$ws.Range("D41").Value = '1.93'
$ws.Range("D42").Value = '4.04'
$ws.Range("D46").Value = '22.35'
$ws.Range("D49").Value = '2.11'
$ws.Range("D51").Value = '0.0362'
